$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert the new header cell A1 = "Category"
$ws.Range("A1").Value = "Category"

# Give A1 the same header formatting as the rest of row 1 (bold, bordered, centered)
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)

# Remove the header-style formatting that used to be on A2:A46 (category column body)
$ws.Range("A2:A46").ClearFormats()
